# Apply the row permutation described by the diff: the data rows 5-17 on
# the single worksheet get their "record" contents (columns A,B,D,E,F,G,H,K,Q,R)
# reshuffled among the row numbers. All other columns (C,I,P,S,T,U,V,W,Y,AA,
# AD,AE,AG,AT,AW,AX,AY) are identical across every one of these rows already,
# so only the moving columns need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move, keyed by row number.
$cols = @("A","B","D","E","F","G","H","K","Q","R")
$rows = 5..17

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# after-row -> before-row (where the content now destined for "after-row"
# used to live before the edit).
$mapping = @{
    5  = 6
    6  = 5
    7  = 11
    8  = 10
    9  = 13
    10 = 9
    11 = 15
    12 = 8
    13 = 16
    14 = 7
    15 = 17
    16 = 12
    17 = 14
}

foreach ($afterRow in $rows) {
    $beforeRow = $mapping[$afterRow]
    $src = $snapshot[$beforeRow]
    foreach ($c in $cols) {
        if ($c -eq "K") {
            # K only has content on one source row ("i frukt"); clear it
            # everywhere else so it only ends up on the row that inherits it.
            if ($src.ContainsKey("K") -and $src["K"] -ne $null -and $src["K"] -ne "") {
                $ws.Range("K$afterRow").Value = $src["K"]
            } else {
                $ws.Range("K$afterRow").Value = ""
            }
        } else {
            $ws.Range("$c$afterRow").Value = $src[$c]
        }
    }
}
